$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131, shifting existing rows 131-147 down to 132-148.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new weekly record.
$ws.Cells.Item(131, 1).Value = 5
$ws.Cells.Item(131, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(131, 3).Value = "Maule"
$ws.Cells.Item(131, 4).Value = 45223
$ws.Cells.Item(131, 5).Value = 7
$ws.Cells.Item(131, 6).Value = 100112022
$ws.Cells.Item(131, 7).Value = "Arveja Verde"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 500
$ws.Cells.Item(131, 11).Value = 20000
$ws.Cells.Item(131, 12).Value = 22000
$ws.Cells.Item(131, 13).Value = 20800
$ws.Cells.Item(131, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(131, 15).Value = "Región del Maule"
$ws.Cells.Item(131, 16).Value = 832
$ws.Cells.Item(131, 17).Value = 25
$ws.Cells.Item(131, 18).Value = "Hortaliza"
